{"js": "// Replace the ten groups of \"three-digit \u00d7 one-digit\" practice problems\n// with the new values from the commit diff. Each original expression is\n// unique in the document, so a direct search + replace per pair is safe\n// and order-independent.\nconst replacements = [\n  [\"518\u00d76=\", \"539\u00d79=\"],\n  [\"767\u00d72=\", \"462\u00d73=\"],\n  [\"931\u00d76=\", \"415\u00d75=\"],\n  [\"632\u00d77=\", \"947\u00d77=\"],\n  [\"624\u00d73=\", \"246\u00d75=\"],\n  [\"218\u00d73=\", \"321\u00d72=\"],\n  [\"702\u00d79=\", \"246\u00d78=\"],\n  [\"826\u00d77=\", \"782\u00d77=\"],\n  [\"278\u00d74=\", \"984\u00d78=\"],\n  [\"476\u00d75=\", \"887\u00d76=\"],\n  [\"932\u00d78=\", \"267\u00d72=\"],\n  [\"723\u00d77=\", \"824\u00d78=\"],\n  [\"209\u00d75=\", \"897\u00d78=\"],\n  [\"461\u00d77=\", \"547\u00d78=\"],\n  [\"549\u00d74=\", \"630\u00d74=\"],\n  [\"829\u00d73=\", \"731\u00d74=\"],\n  [\"495\u00d75=\", \"981\u00d74=\"],\n  [\"482\u00d75=\", \"682\u00d76=\"],\n  [\"139\u00d76=\", \"210\u00d75=\"],\n  [\"231\u00d77=\", \"249\u00d76=\"],\n  [\"322\u00d74=\", \"839\u00d75=\"],\n  [\"329\u00d77=\", \"413\u00d77=\"],\n  [\"916\u00d74=\", \"209\u00d73=\"],\n  [\"465\u00d75=\", \"858\u00d76=\"],\n  [\"446\u00d77=\", \"714\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the ten groups of \"three-digit x one-digit\" practice problems\n# with the new values from the commit diff. Each original expression is\n# unique in the document, so a direct Find/Replace per pair is safe and\n# order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"518\u00d76=\", \"539\u00d79=\"),\n    @(\"767\u00d72=\", \"462\u00d73=\"),\n    @(\"931\u00d76=\", \"415\u00d75=\"),\n    @(\"632\u00d77=\", \"947\u00d77=\"),\n    @(\"624\u00d73=\", \"246\u00d75=\"),\n    @(\"218\u00d73=\", \"321\u00d72=\"),\n    @(\"702\u00d79=\", \"246\u00d78=\"),\n    @(\"826\u00d77=\", \"782\u00d77=\"),\n    @(\"278\u00d74=\", \"984\u00d78=\"),\n    @(\"476\u00d75=\", \"887\u00d76=\"),\n    @(\"932\u00d78=\", \"267\u00d72=\"),\n    @(\"723\u00d77=\", \"824\u00d78=\"),\n    @(\"209\u00d75=\", \"897\u00d78=\"),\n    @(\"461\u00d77=\", \"547\u00d78=\"),\n    @(\"549\u00d74=\", \"630\u00d74=\"),\n    @(\"829\u00d73=\", \"731\u00d74=\"),\n    @(\"495\u00d75=\", \"981\u00d74=\"),\n    @(\"482\u00d75=\", \"682\u00d76=\"),\n    @(\"139\u00d76=\", \"210\u00d75=\"),\n    @(\"231\u00d77=\", \"249\u00d76=\"),\n    @(\"322\u00d74=\", \"839\u00d75=\"),\n    @(\"329\u00d77=\", \"413\u00d77=\"),\n    @(\"916\u00d74=\", \"209\u00d73=\"),\n    @(\"465\u00d75=\", \"858\u00d76=\"),\n    @(\"446\u00d77=\", \"714\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $true, $true, $false, $null, $null, $true, $null, $null, $null, 2)\n}\n"}
